$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shorten the use-case title
$ws.Range("C2").Value = "Selecionar Receita"

# Renumber the second exception label
$ws.Range("B19").Value = "Exceção 2"
